$wb = $excel.ActiveWorkbook

# --- CoordinatedOps: append another year (1/1/2021) of data in row 5 ---
$ws = $wb.Worksheets.Item("CoordinatedOps")
$ws.Cells.Item(5,1).Value = 44197
$ws.Range("A5").NumberFormat = "m/d/yyyy\ h:mm:ss"
$ws.Cells.Item(5,2).Value = "NaN"
$ws.Cells.Item(5,3).Value = "NaN"
$ws.Cells.Item(5,4).Value = "NaN"
$ws.Cells.Item(5,5).Value = "NaN"
$ws.Cells.Item(5,6).Value = 3427
$ws.Cells.Item(5,7).Value = 955

# --- Update selection on Reservoirs sheet (bottom-right pane) ---
$wsReservoirs = $wb.Worksheets.Item("Reservoirs")
$wsReservoirs.Activate()
$wsReservoirs.Range("L18").Select()

# --- Make CoordinatedOps the active/selected sheet with its own selection ---
$ws.Activate()
$ws.Range("F6").Select()
